$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.180.76'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.333.92'
$ws.Range("E3").Value = '  +3.87%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '109.44'
$ws.Range("E5").Value = '  -4.17%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '310.55'
$ws.Range("E6").Value = '  +3.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.636'
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("E8").Value = '  -0.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.05'
$ws.Range("E10").Value = '  -3.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0932'
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.90'
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.07'
$ws.Range("E13").Value = '  +19.34%  '
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.89'
$ws.Range("E15").Value = '  +3.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.672.55'
$ws.Range("E16").Value = '  +3.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.408.21'
$ws.Range("E17").Value = '  +6.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.154.78'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("E19").Value = '  -0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.32'
$ws.Range("E20").Value = '  -5.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.01'
$ws.Range("E21").Value = '  +3.47%  '
$ws.Range("E22").Value = '  -5.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.52'
$ws.Range("E23").Value = '  +7.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '254.85'
$ws.Range("E24").Value = '  +9.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.08'
$ws.Range("E25").Value = '  -4.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.91'
$ws.Range("E26").Value = '  -2.06%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '39.26'
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("E29").Value = '  +1.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.61'
$ws.Range("E30").Value = '  +6.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.01'
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("E32").Value = '  -3.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0909'
$ws.Range("E33").Value = '  +0.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.81'
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.131'
$ws.Range("E36").Value = '  +1.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.17'
$ws.Range("E37").Value = '  -6.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0378'
$ws.Range("E38").Value = '  +1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.105'
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.75'
$ws.Range("E40").Value = '  +5.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.50'
$ws.Range("E41").Value = '  +11.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.95'
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.234'
$ws.Range("E43").Value = '  -2.43%  '
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.51'
$ws.Range("E45").Value = '  -6.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.74'
$ws.Range("E46").Value = '  +2.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '110.97'
$ws.Range("E47").Value = '  +4.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.15'
$ws.Range("E48").Value = '  +4.90%  '
$ws.Range("E49").Value = '  -6.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0993'
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.07'
$ws.Range("E51").Value = '  +2.83%  '
